$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.912.12"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.99"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.95"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4624"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07318"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8741"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07831"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +3.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.59"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.811.83"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.323"
$ws.Range("D14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.537"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.03"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.27%  "

$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008835"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.72"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.935.39"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.096"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.51"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.037.03"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.71"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.030"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.087"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.12"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08843"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.62%  "

$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7303"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.424"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.459"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.072"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01936"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.953"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.069"
$ws.Range("D41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5116"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1620"
$ws.Range("D43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.128"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4814"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.19"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.616"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06191"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.44"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.34%  "
